# ParseOptions: make sheetName selectable - add a second sheet ("시트러스")
# with extended product data, and preserve the original "Sheet1" data
# (unchanged) as a second sheet in the workbook so the parser can be
# exercised against either sheet name.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# 1) Duplicate the ORIGINAL (pre-edit) data of the existing sheet into a
#    brand-new sheet placed right after it, before touching $ws1's data.
$ws2 = $wb.Worksheets.Add($null, $ws1)

# Rename the original/first sheet first so the new sheet can take over
# the "Sheet1" name without a naming collision.
$ws1.Name = "시트러스"

# Column A carries a custom number-format style (s="1") in the source -
# copy formats only for that column so the new sheet reuses the same
# style index instead of minting a duplicate one.
$ws1.Range("A1:A3").Copy()
$ws2.Range("A1").PasteSpecial(-4122)

$ws2.Range("A1").Value = "상품ID"
$ws2.Range("B1").Value = "상품명칭"
$ws2.Range("C1").Value = "바코드"
$ws2.Range("D1").Value = "가격"
$ws2.Range("A2").Value = 1001
$ws2.Range("B2").Value = "가나다라"
$ws2.Range("C2").Value = "B12345123"
$ws2.Range("D2").Value = 1200
$ws2.Range("A3").Value = 1002
$ws2.Range("C3").Value = "B23345123"
$ws2.Range("D3").Value = 1200

$ws2.Range("A1:D3").Select()
$ws2.Name = "Sheet1"

# 2) Append the new rows to the renamed original sheet.

$ws1.Range("B3").Value = "마바사아"
$ws1.Range("B4").Value = "자차카타"
$ws1.Range("B5").Value = "파하"

$ws1.Range("C4").Value = "B12345124"
$ws1.Range("C5").Value = "B23345125"

$ws1.Range("D3").Value = 1600
$ws1.Range("A4").Value = 1003
$ws1.Range("D4").Value = 1800
$ws1.Range("A5").Value = 1004
$ws1.Range("D5").Value = 12000

# 3) Restore the first sheet as the active tab/selection, matching the
#    final cursor position recorded in the workbook.
$ws1.Select()
$ws1.Range("F24").Select()
